$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "42.477.80"
$ws.Range("E2").Value = "  +0.10%  "

Set-TextCell $ws "D3" "2.297.96"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextCell $ws "D5" "315.85"
$ws.Range("E5").Value = "  -1.02%  "

Set-TextCell $ws "D6" "102.86"
$ws.Range("E6").Value = "  -1.36%  "

Set-TextCell $ws "D7" "0.628"
$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("E8").Value = "  -0.02%  "

Set-TextCell $ws "D9" "0.603"
$ws.Range("E9").Value = "  -1.64%  "

Set-TextCell $ws "D10" "39.40"
$ws.Range("E10").Value = "  -1.73%  "

Set-TextCell $ws "D11" "0.0905"
$ws.Range("E11").Value = "  -0.79%  "

Set-TextCell $ws "D12" "8.48"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("E14").Value = "  +2.47%  "

Set-TextCell $ws "D15" "15.27"
$ws.Range("E15").Value = "  -1.18%  "

Set-TextCell $ws "D16" "2.647.32"
$ws.Range("E16").Value = "  -0.66%  "

Set-TextCell $ws "D17" "2.295.87"
$ws.Range("E17").Value = "  -1.06%  "

Set-TextCell $ws "D18" "42.560.39"
$ws.Range("E18").Value = "  +0.14%  "

Set-TextCell $ws "D19" "7.55"
$ws.Range("E19").Value = "  +0.49%  "

Set-TextCell $ws "D20" "13.87"
$ws.Range("E20").Value = "  +24.06%  "

Set-TextCell $ws "D21" "0.0000105"
$ws.Range("E21").Value = "  -1.02%  "

Set-TextCell $ws "D22" "73.95"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("E23").Value = "  -3.78%  "

Set-TextCell $ws "D24" "265.74"
$ws.Range("E24").Value = "  -5.52%  "

$ws.Range("E25").Value = "  -2.21%  "

$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  -4.45%  "

Set-TextCell $ws "D29" "22.60"
$ws.Range("E29").Value = "  -2.16%  "

Set-TextCell $ws "D30" "6.61"
$ws.Range("E30").Value = "  +11.33%  "

Set-TextCell $ws "D31" "37.20"
$ws.Range("E31").Value = "  +1.68%  "

Set-TextCell $ws "D32" "165.59"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +0.02%  "

Set-TextCell $ws "D34" "0.131"
$ws.Range("E34").Value = "  -3.76%  "

Set-TextCell $ws "D35" "2.60"
$ws.Range("E35").Value = "  -4.25%  "

$ws.Range("E36").Value = "  -2.05%  "

Set-TextCell $ws "D37" "4.55"
$ws.Range("E37").Value = "  -2.46%  "

Set-TextCell $ws "D38" "0.0353"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("E39").Value = "  -0.66%  "

Set-TextCell $ws "D40" "2.70"
$ws.Range("E40").Value = "  -3.82%  "

$ws.Range("E41").Value = "  +6.25%  "

Set-TextCell $ws "D42" "70.25"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("E45").Value = "  +0.22%  "

Set-TextCell $ws "D46" "12.29"
$ws.Range("E46").Value = "  +0.87%  "

Set-TextCell $ws "D47" "115.02"
$ws.Range("E47").Value = "  +1.26%  "

Set-TextCell $ws "D48" "80.23"
$ws.Range("E48").Value = "  +1.89%  "

Set-TextCell $ws "D49" "1.657.89"
$ws.Range("E49").Value = "  +2.56%  "

Set-TextCell $ws "D50" "5.25"
$ws.Range("E50").Value = "  -1.44%  "

Set-TextCell $ws "D51" "8.81"
$ws.Range("E51").Value = "  -2.17%  "

# Row 43/44 swap: Algorand moves to row 43, BitcoinSV moves to row 44
Set-TextCell $ws "B43" "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws "D43" "0.228"
$ws.Range("E43").Value = "  +0.13%  "

Set-TextCell $ws "B44" "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell $ws "D44" "94.62"
$ws.Range("E44").Value = "  -5.73%  "
